$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 141 (pushes the existing rows 141..198 down to 142..199)
$ws.Rows("141:141").Insert()

# Fill the newly inserted row 141 with the new weekly price record
$ws.Range("A141").Value = 8
$ws.Range("B141").Value = "Terminal La Palmera de La Serena"
$ws.Range("C141").Value = "Coquimbo"
$ws.Range("D141").Value = 44608
$ws.Range("E141").Value = 4
$ws.Range("F141").Value = 100112037
$ws.Range("G141").Value = "Cebollín"
$ws.Range("H141").Value = "Sin especificar"
$ws.Range("I141").Value = "Primera"
$ws.Range("J141").Value = 200
$ws.Range("K141").Value = 8500
$ws.Range("L141").Value = 9000
$ws.Range("M141").Value = 8750
$ws.Range("N141").Value = "$/paquete 36 unidades"
$ws.Range("O141").Value = "Provincia del Elquí"
$ws.Range("P141").Value = 243
$ws.Range("Q141").Value = 36
$ws.Range("R141").Value = "Hortaliza"
